$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "urban"/"rural" row labels with the more specific
# "city/urban settlement" and "village/rural area" labels (Kyrgyz, Russian, English).
$ws.Range("A6").Value = "Шаар жерлери"
$ws.Range("B6").Value = "Городские поселения"
$ws.Range("C6").Value = "City"
$ws.Range("A7").Value = "Айыл аймагы"
$ws.Range("B7").Value = "Сельская местность"
$ws.Range("C7").Value = "Village"

# Move the selection/scroll position: no more frozen/scrolled topLeftCell,
# and the active selection moves to B29.
[void]$ws.Range("B29").Select()
